$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("2")
$before = $wb.Worksheets.Item("2")

$source.Copy($before)

$newSheet = $wb.Worksheets.Item("2 (2)")
$newSheet.Name = "GroupWithDifferentNoteNumbers"

$newSheet.Range("B3").Value = 1663
$newSheet.Select()
$newSheet.Range("B3").Select()
